$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 45171 to 45172
# for all data rows (rows 2 through 269).
$ws.Range("C2:C269").Value = 45172
